$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.113.28"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "3.784.41"
$ws.Range("E3").Value = "  -2.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.72"
$ws.Range("E5").Value = "  -1.02%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.47"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "3.783.85"
$ws.Range("E7").Value = "  -2.44%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  -2.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  -3.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.75"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("D15").Value = "4.416.44"
$ws.Range("E15").Value = "  -2.39%  "
$ws.Range("D16").Value = "3.772.55"
$ws.Range("E16").Value = "  -2.61%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.66"
$ws.Range("E17").Value = "  +1.67%  "
$ws.Range("D18").Value = "68.021.09"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  -3.72%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -4.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.36"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.718"
$ws.Range("E23").Value = "  -2.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000147"
$ws.Range("E24").Value = "  -9.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.87"
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.17"
$ws.Range("E27").Value = "  -0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.41"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.92"
$ws.Range("E30").Value = "  -1.15%  "
$ws.Range("D31").Value = "3.929.47"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.45"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("E34").Value = "  -4.28%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.21"
$ws.Range("E35").Value = "  -2.45%  "
$ws.Range("D36").Value = "3.739.46"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.71"
$ws.Range("E38").Value = "  -6.44%  "
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "405.50"
$ws.Range("E47").Value = "  -3.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.56"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.81"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000271"
$ws.Range("E50").Value = "  -11.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.97"
$ws.Range("E51").Value = "  +3.23%  "
